$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.278.87'
$ws.Range('E2').Value = '  -0.46%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.860.27'
$ws.Range('E3').Value = '  -1.06%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.44'
$ws.Range('E5').Value = '  -0.47%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6978'
$ws.Range('E6').Value = '  -2.71%  '
$ws.Range('E7').Value = '  -0.18%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07842'
$ws.Range('E8').Value = '  -1.17%  '
$ws.Range('E9').Value = '  -1.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.98'
$ws.Range('E10').Value = '  -3.97%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07791'
$ws.Range('E11').Value = '  -4.27%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.866.07'
$ws.Range('E12').Value = '  -1.29%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.134'
$ws.Range('E13').Value = '  -2.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '92.28'
$ws.Range('E14').Value = '  -2.54%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6922'
$ws.Range('E15').Value = '  -2.44%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.524'
$ws.Range('E16').Value = '  +1.92%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008473'
$ws.Range('E17').Value = '  +0.56%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.213.75'
$ws.Range('E18').Value = '  -0.76%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '248.41'
$ws.Range('E19').Value = '  -2.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.097.85'
$ws.Range('E20').Value = '  -2.62%  '
$ws.Range('E21').Value = '  -2.91%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9996'
$ws.Range('E22').Value = '  -0.15%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.533'
$ws.Range('E23').Value = '  -2.91%  '
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1532'
$ws.Range('E25').Value = '  -3.46%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '161.41'
$ws.Range('E26').Value = '  -0.84%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.898'
$ws.Range('E27').Value = '  -2.00%  '
$ws.Range('E28').Value = '  -1.63%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.269'
$ws.Range('E30').Value = '  -3.42%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.243'
$ws.Range('E31').Value = '  -0.94%  '
$ws.Range('E32').Value = '  -2.46%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05221'
$ws.Range('E33').Value = '  -2.19%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.874'
$ws.Range('E34').Value = '  -3.95%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7518'
$ws.Range('E35').Value = '  -0.73%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.172'
$ws.Range('E36').Value = '  -0.73%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.695'
$ws.Range('E37').Value = '  -0.24%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01862'
$ws.Range('E38').Value = '  -1.84%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.245.28'
$ws.Range('E39').Value = '  -2.44%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.745'
$ws.Range('E40').Value = '  -0.89%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9003'
$ws.Range('E41').Value = '  -0.43%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '111.32'
$ws.Range('E42').Value = '  -1.30%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.932'
$ws.Range('E43').Value = '  -8.23%  '
$ws.Range('E44').Value = '  -0.10%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '68.79'
$ws.Range('E45').Value = '  -7.46%  '
$ws.Range('E46').Value = '  -5.14%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.999.52'
$ws.Range('E47').Value = '  -1.96%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.565'
$ws.Range('E48').Value = '  +0.75%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.5183'
$ws.Range('E49').Value = '  -0.51%  '
$ws.Range('E50').Value = '  -1.48%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4258'
$ws.Range('E51').Value = '  -2.61%  '

Write-Host "Applied 92 cell updates"